$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their original text representation
# (these values would otherwise be auto-converted to numbers by Excel)
$textCells = @("D4","D5","D6","D9","D11","D13","D14","D17","D19","D21","D22","D24","D25","D27","D28","D32","D34","D37","D38","D39","D47","D48","D49","D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "67.785.55"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.790.90"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "597.62"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "168.62"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").Value = "3.790.19"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +5.91%  "
$ws.Range("D14").Value = "36.81"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "4.423.76"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "3.795.50"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "19.01"
$ws.Range("E17").Value = "  +5.09%  "
$ws.Range("D18").Value = "67.656.19"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "7.30"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "468.28"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "0.0000150"
$ws.Range("E24").Value = "  -6.00%  "
$ws.Range("D25").Value = "83.40"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").Value = "3.942.67"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").Value = "30.55"
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("D36").Value = "3.750.21"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.78"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "5.91"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "5.91"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "407.59"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("D48").Value = "46.29"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "0.000279"
$ws.Range("E49").Value = "  -7.17%  "
$ws.Range("D50").Value = "141.96"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("E51").Value = "  -0.41%  "
